$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (3-7): Name, Last name(blank), Email, Phone, Address, City(blank), Created at
$rowsData = @(
    @("asdasdasdasd", "asdasd", "asdasd", "asdasd", "12/06/2014 - 03:05:19"),
    @("asdasd", "asdasd", "asdasd", "asdasd", "12/06/2014 - 03:05:59"),
    @("asdasd", "asdasd", "asdasd", "asdasda", "12/06/2014 - 03:06:23"),
    @("Josesito", "Jiseselito@josnrn.hhsj", "62727:73", "Jesús jdjd k 123", "13/06/2014 - 22:23:41"),
    @("evangelina tapia", "eva@chile.com", 61231234, "Los Conejitos Blanquitos 123", "16/06/2014 - 21:19:54")
)

$startRow = 3
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rowsData[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = $rowData[1]
    $ws.Cells.Item($r, 4).Value = $rowData[2]
    $ws.Cells.Item($r, 5).Value = $rowData[3]
    $ws.Cells.Item($r, 6).Style = "Normal"
    $ws.Cells.Item($r, 7).Value = $rowData[4]
}

# Column width adjustments (closest achievable values through this engine's
# character-width quantization) for columns C and E.
$ws.Columns.Item(3).ColumnWidth = 12.25
$ws.Columns.Item(5).ColumnWidth = 21.15
